$d = $word.ActiveDocument

# --- Text correction 1: extend "first contact" sentence ---
$old1 = "An introduction for participants who are completely new to R which covers finding their way round RStudio importing some data, summarising and plotting it. It provides a first contact"
$new1 = "An introduction for participants who are completely new to R which covers finding their way round RStudio importing some data, summarising and plotting it. It provides a first contact with ideas explored in more detail in later modules such as RStudio Projects, data import, tidy data and data reformatting."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# --- Text correction 2: extend "cause a lot of stress" sentence ---
$old2 = "About Project-oriented workflow, working directories and paths, project organisation and naming things! Over many years in teaching computational biology I have seen many people struggle not because of the analysis itself but because the ideas of working directory and paths are unfamiliar. This can cause a lot of stress"
$new2 = "About Project-oriented workflow, working directories and paths, project organisation and naming things! Over many years in teaching computational biology I have seen many people struggle not because of the analysis itself but because the ideas of working directory and paths are unfamiliar. This can cause a lot of stress and gives people the impression the computational analysis is too hard for them when they just have a small knowledge gap. Plugging this gap is very empowering!"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# --- Resize the four figures from ~5.017"x4.014" (361.3pt x 289pt) down to a clean 5"x4" (360pt x 288pt) ---
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    if ([Math]::Abs($shp.Width - 361.3) -lt 0.5 -and [Math]::Abs($shp.Height - 289) -lt 0.5) {
        $shp.Width = 360
        $shp.Height = 288
    }
}
